# Insert a new data row at row 294 (pushing existing rows 294-333 down to 295-334)
# and populate it with the new observation described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 294..333 down by one, inserting a fresh (blank) row at 294.
$ws.Rows.Item(294).Insert()

# Populate the newly inserted row 294 with the new record.
$ws.Range("A294").Value = 10
$ws.Range("B294").Value = "Vega Modelo de Temuco"
$ws.Range("C294").Value = "La Araucanía"
$ws.Range("D294").Value = 44505
$ws.Range("E294").Value = 9
$ws.Range("F294").Value = 100112028
$ws.Range("G294").Value = "Sandia"
$ws.Range("H294").Value = "Sin especificar"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 450
$ws.Range("K294").Value = 800
$ws.Range("L294").Value = 900
$ws.Range("M294").Value = 856
$ws.Range("N294").Value = "$/kilo (volumen en unidades)"
$ws.Range("O294").Value = "Perú"
$ws.Range("P294").Value = 856
$ws.Range("Q294").Value = 1
$ws.Range("R294").Value = "Hortaliza"

"done"
